# Consolidate the three text runs that make up the Slide 2 title
# ("Lists" + " " + "(continued)") into a single run "Lists (continued)".
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(1)

# The writer keeps existing runs untouched wherever the new text shares a
# common prefix/suffix with the current (concatenated) text, which would
# just re-trim the surrounding runs instead of merging them. Assigning an
# unrelated placeholder value first forces a full rebuild of the text
# body as a single run; the subsequent assignment of the real text then
# lands in that single run instead of being split across the old ones.
$shape.TextFrame.TextRange.Text = "placeholder-reset-zzz"
$shape.TextFrame.TextRange.Text = "Lists (continued)"
